$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-02-07", "sleep",           $false, $false),
    @("2025-02-07", "activity",        $true,  $true),
    @("2025-02-07", "weekly_activity", $true,  $false),
    @("2025-02-08", "sleep",           $false, $false),
    @("2025-02-08", "activity",        $false, $false),
    @("2025-02-08", "weekly_activity", $false, $false)
)

$startRow = 20
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
